# Update gh-pages to output generated at 456a3b4
# Applies numeric "interested count" (F) / price (G) refreshes across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets, a
# time-range correction for the Phantom of the Opera listing, and appends a
# brand-new 演出 row for the Shikawa Ayako violin concert.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (exhibitions) -------------------------------------------
$ws1.Range("F2").Value2  = 2951
$ws1.Range("F3").Value2  = 6411
$ws1.Range("F4").Value2  = 2520
$ws1.Range("F6").Value2  = 528
$ws1.Range("F7").Value2  = 55
$ws1.Range("F9").Value2  = 2927
$ws1.Range("F10").Value2 = 352
$ws1.Range("F12").Value2 = 7354
$ws1.Range("F13").Value2 = 336
$ws1.Range("F16").Value2 = 244
$ws1.Range("F19").Value2 = 8973
$ws1.Range("F27").Value2 = 105
$ws1.Range("F30").Value2 = 48
$ws1.Range("F32").Value2 = 108
$ws1.Range("F33").Value2 = 2616
$ws1.Range("F36").Value2 = 42

# Phantom of the Opera listing: showtime end pushed 17:00 -> 21:00, and the
# minimum ticket price drops from 55 to 30.
$ws1.Range("E37").Value2 = "2024.05.25 09:00-05.26 21:00"
$ws1.Range("G37").Value2 = 30

$ws1.Range("F38").Value2 = 748
$ws1.Range("F39").Value2 = 3871
$ws1.Range("F41").Value2 = 31
$ws1.Range("F42").Value2 = 1196
$ws1.Range("F43").Value2 = 60
$ws1.Range("F44").Value2 = 18
$ws1.Range("F45").Value2 = 228
$ws1.Range("F46").Value2 = 8
$ws1.Range("F47").Value2 = 51
$ws1.Range("F49").Value2 = 49

# --- Sheet "演出" (performances) -------------------------------------------
$ws2.Range("F7").Value2 = 122

# New row 19: 北京·石川绫子小提琴动漫音乐会
$ws2.Range("A18").Copy()
$ws2.Range("A19").PasteSpecial(-4122)

$ws2.Range("A19").Value2 = 18
$ws2.Range("B19").NumberFormat = "@"
$ws2.Range("B19").Value2 = "2024-07-22"
$ws2.Range("C19").Value2 = "北京·石川绫子小提琴动漫音乐会"
$ws2.Range("D19").Value2 = "中关村南大街33号中国国家图书馆内 国图艺术中心"
$ws2.Range("E19").Value2 = "2024.07.22 19:30-07.22 21:00"
$ws2.Range("F19").Value2 = 0
$ws2.Range("G19").Value2 = 180
$ws2.Range("H19").Value2 = "https://show.bilibili.com/platform/detail.html?id=83973"
$ws2.Range("I19").Value2 = "//i1.hdslb.com/bfs/openplatform/202404/HhY3CS7t1712652128640.jpeg"

# --- Sheet "全部类型" (all types) ------------------------------------------
$ws4.Range("F3").Value2  = 2951
$ws4.Range("F6").Value2  = 6411
$ws4.Range("F7").Value2  = 2520
$ws4.Range("F8").Value2  = 122
$ws4.Range("F10").Value2 = 528
$ws4.Range("F11").Value2 = 55
$ws4.Range("F13").Value2 = 2927
$ws4.Range("F14").Value2 = 352
$ws4.Range("F18").Value2 = 7354
$ws4.Range("F19").Value2 = 336
$ws4.Range("F22").Value2 = 244
$ws4.Range("F24").Value2 = 8973
$ws4.Range("F30").Value2 = 105
$ws4.Range("F31").Value2 = 48
$ws4.Range("F33").Value2 = 108
$ws4.Range("F34").Value2 = 2616
$ws4.Range("F37").Value2 = 42

$ws4.Range("E38").Value2 = "2024.05.25 09:00-05.26 21:00"
$ws4.Range("G38").Value2 = 30

$ws4.Range("F39").Value2 = 748
$ws4.Range("F41").Value2 = 3871
$ws4.Range("F43").Value2 = 31
$ws4.Range("F45").Value2 = 1196
$ws4.Range("F46").Value2 = 228
$ws4.Range("F47").Value2 = 51
$ws4.Range("F49").Value2 = 49
